$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting the existing rows 86..109 down to 87..110.
$ws.Rows(86).Insert()

# Populate the newly inserted row 86 with the new weekly price-report entry.
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C86").Value = "Ñuble"
$ws.Range("D86").Value = 44924
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = 100112031
$ws.Range("G86").Value = "Poroto verde"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 60
$ws.Range("K86").Value = 25000
$ws.Range("L86").Value = 26000
$ws.Range("M86").Value = 25500
$ws.Range("N86").Value = "$/saco 25 kilos"
$ws.Range("O86").Value = "Región del Maule"
$ws.Range("P86").Value = 1020
$ws.Range("Q86").Value = 25
$ws.Range("R86").Value = "Hortaliza"
